# Update the "Data" sheet of the USDSOFRCSA_USD workbook.
# The 1Y SOFROIS/OIS row is removed and replaced with a block of new
# SOFR FUTURE quotes (11 rows), pushing the remaining OIS tenors
# (2Y, 3Y, 5Y, 10Y, 20Y, 30Y) further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Remove the old "1Y" row (row 3).
$ws.Rows("3:3").Delete()

# Insert 11 blank rows where the FUTURE block will go (new rows 3-13),
# pushing the rest of the OIS data (now starting at row 3) back down.
$ws.Rows("3:13").Insert()

# New FUTURE rows (Tenor, Ticker, Type, Rate).
$futures = @(
    @("3M",  "SQZ25", "FUTURE", 96.19499999999999),
    @("5M",  "SQF26", "FUTURE", 96.29000000000001),
    @("6M",  "SQG26", "FUTURE", 96.38),
    @("7M",  "SQH26", "FUTURE", 96.41),
    @("9M",  "SQM26", "FUTURE", 96.63),
    @("12M", "SQU26", "FUTURE", 96.785),
    @("15M", "SQZ26", "FUTURE", 96.86499999999999),
    @("0M",  "SQQ25", "FUTURE", 95.78),
    @("0M",  "SQU25", "FUTURE", 95.8925),
    @("1M",  "SQV25", "FUTURE", 96),
    @("3M",  "SQX25", "FUTURE", 96.11499999999999)
)

$startRow = 3
for ($i = 0; $i -lt $futures.Count; $i++) {
    $row = $startRow + $i
    $data = $futures[$i]
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
}
